# Implemented better battle interrupt system:
# INTERRUPT rows used to store a fixed health-target keyword (MIDDLE_HEALTH /
# RIGHT_HEALTH / PLAYER_HEALTH) in column B, a numeric threshold in column C,
# and the "who must be alive" flag in column D. They now store a single
# arbitrary interrupt-condition string (handled by BattleInterrupts.cs) in
# column B, and the "who must be alive" flag moves up into column C - column D
# is no longer used for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: MIDDLE_HEALTH / 0.5 / L  ->  check-health,1,0.5 / L
$ws.Range("B10").Value = "check-health,1,0.5"
$ws.Range("C10").Value = "L"
$ws.Range("D10").ClearContents()

# Row 13: MIDDLE_HEALTH / 0 / R  ->  check-health,1,0 / R
$ws.Range("B13").Value = "check-health,1,0"
$ws.Range("C13").Value = "R"
$ws.Range("D13").ClearContents()

# Row 17: RIGHT_HEALTH / 0 / R  ->  check-health,2,0 / R
$ws.Range("B17").Value = "check-health,2,0"
$ws.Range("C17").Value = "R"
$ws.Range("D17").ClearContents()

# Row 20: PLAYER_HEALTH / 0.75 / LR  ->  check-health,3,0.75 / LR
$ws.Range("B20").Value = "check-health,3,0.75"
$ws.Range("C20").Value = "LR"
$ws.Range("D20").ClearContents()

# Update the view state: scrolled down a bit with B24 as the active selection.
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
